# Apply the commit "feat: Improve Final Project Instruction" edits.
$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: expand the grading-scale sentence into 5 runs; swap
# "صفر، نیم یا یک" -> "عددی بین صفر تا یک" and
# "عدد نیم یا صفر" -> "عددی کمتر از یک".
# ---------------------------------------------------------------------
$f1 = $d.Content
$found1 = $f1.Find.Execute("نمرات هر قسمت از پروژه صفر، نیم یا یک خواهد بود. یعنی در اکسل برای هر قسمت یا آن قسمت کار می‌کند که عدد یک وارد می‌کنید یا کار نمی‌کند که با توجه به برداشت خودتان عدد نیم یا صفر را وارد می‌کنید.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) { throw "Change 1: original sentence not found" }
$target1 = $d.Range($f1.Start, $f1.End)
$xml1 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r w:rsidRPr="00B4246F"><w:rPr><w:rtl/></w:rPr><w:t xml:space="preserve">نمرات هر قسمت از پروژه </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="cs"/><w:rtl/></w:rPr><w:t>عددی بین صفر تا یک</w:t></w:r><w:r><w:rPr><w:rtl/></w:rPr><w:t xml:space="preserve"> خواهد بود. یعنی در اکسل برای هر قسمت یا آن قسمت کار می‌کند که عدد یک وارد می‌کنید یا کار نمی‌کند که با توجه به برداشت خودتان </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="cs"/><w:rtl/></w:rPr><w:t>عددی کمتر از یک</w:t></w:r><w:r><w:rPr><w:rtl/></w:rPr><w:t xml:space="preserve"> را وارد می‌کنید.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$null = $target1.InsertXML($xml1)

# ---------------------------------------------------------------------
# Change 2 & 3: drop the stray <w:rFonts w:hint="cs"/> from the
# paragraph-mark rPr of the "در بازی چند نفره..." paragraph, and add a
# new paragraph right after it describing the end-game / linked-list
# implementation notes (with the _GoBack bookmark trailing it).
# ---------------------------------------------------------------------
$f2 = $d.Content
$found2 = $f2.Find.Execute("در بازی چند نفره روند بازی تک نفره می‌بایست کاملا صحیح رعایت شده و نوبت‌ها به درستی به بازیکن‌ها تخصیص داده شود.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) { throw "Change 2: target paragraph not found" }
$targetPara = $f2.Paragraphs.Item(1)
$target2 = $d.Range($targetPara.Range.Start, $targetPara.Range.End)
$xml2 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p w14:paraId="5F2BA19A" w14:textId="6E967C1A" w:rsidR="003C5EED" w:rsidRPr="00E34D1C" w:rsidRDefault="003C5EED" w:rsidP="003C5EED">
<w:pPr><w:bidi/><w:rPr><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr></w:pPr>
<w:r><w:rPr><w:rFonts w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>در بازی چند نفره روند بازی تک نفره می‌بایست کاملا صحیح رعایت شده و نوبت‌ها به درستی به بازیکن‌ها تخصیص داده شود.</w:t></w:r>
</w:p>
<w:p>
<w:pPr><w:bidi/><w:rPr><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr></w:pPr>
<w:r><w:rPr><w:rFonts w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve">در بازی برای پایان، شرط اصلی همان گزینه خروج می‌باشد. برای تسلط دانشجویان </w:t></w:r><w:r><w:rPr><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>پ</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>ی</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>اده</w:t></w:r><w:r><w:rPr><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> ساز</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>ی</w:t></w:r><w:r><w:rPr><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:lang w:bidi="fa-IR"/></w:rPr><w:t>linked list</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>،</w:t></w:r><w:r><w:rPr><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>ذخیره و بازیابی</w:t></w:r><w:r><w:rPr><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> باز</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>ی</w:t></w:r><w:r><w:rPr><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> و انجام تقس</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>ی</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>م</w:t></w:r><w:r><w:rPr><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t xml:space="preserve"> سلول </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="cs"/><w:rtl/><w:lang w:bidi="fa-IR"/></w:rPr><w:t>به نظر مناسب می‌باشند.</w:t></w:r>
<w:bookmarkStart w:id="2" w:name="_GoBack"/><w:bookmarkEnd w:id="2"/>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$null = $target2.InsertXML($xml2)

Write-Output "edits applied"
